# Add a new lookup row ("clNoL2416Eq") to the "DBS" sheet, mirroring the
# existing clNoFirst / clNoEq rows: Key ID, "other ORDER condition",
# and "read Key condition" (reusing the ClCode1 ASC... Order text).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBS")

# Set column B before column A so that new shared-string entries are
# appended in the same order as in the target workbook.
$ws.Range("B4").Value = "ClCode1 = ,AND ClCode2 = ,AND ClNo = "
$ws.Range("A4").Value = "clNoL2416Eq"
$ws.Range("C4").Value = "ClCode1 ASC,ClCode2 ASC,ClNo ASC,LandSeq ASC"

$ws.Range("A4").Select()
